# "affichage graphique dans Excel"
#
# 1) Two labels ("Contenu du stage" and "Type entreprise") were typed one
#    column too far to the left (column C instead of column B, which is
#    where all the other section headers - "Lieu du stage", etc. - live).
#    Move them over to column B.
# 2) Add the three pie charts that visualise the three data blocks
#    (Lieu du stage / Contenu du stage / Type du stage) and hook them to
#    the worksheet via a drawing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the misplaced section headers -----------------------------------

# "Contenu du stage" : C16 -> B16
$ws.Range("B16").Value = $ws.Range("C16").Value2
$ws.Range("C16").ClearContents()

# "Type entreprise" : C25 -> B25
$ws.Range("B25").Value = $ws.Range("C25").Value2
$ws.Range("C25").ClearContents()

# --- Chart 1 : "Lieu du stage" (D10:E14, name in B10) ---------------------

$chartObj1 = $ws.ChartObjects().Add(205, 15, 300, 230)
$chart1 = $chartObj1.Chart
$chart1.ChartType = 5
$chart1.SetSourceData($ws.Range("D10:E14"))
$chart1.SeriesCollection(1).Name = "=Worksheet!`$B`$10"
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Lieu du stage"
$chart1.HasLegend = $true
$chart1.Legend.Position = -4152
$chart1.SeriesCollection(1).HasDataLabels = $true
$dLbls1 = $chart1.SeriesCollection(1).DataLabels()
$dLbls1.ShowValue = $true
$dLbls1.ShowPercentage = $true
$dLbls1.ShowCategoryName = $false
$dLbls1.ShowSeriesName = $false
$dLbls1.ShowLegendKey = $false

# --- Chart 2 : "Contenu du stage" (D16:E23, name in B16) ------------------

$chartObj2 = $ws.ChartObjects().Add(205, 260, 300, 230)
$chart2 = $chartObj2.Chart
$chart2.ChartType = 5
$chart2.SetSourceData($ws.Range("D16:E23"))
$chart2.SeriesCollection(1).Name = "=Worksheet!`$B`$16"
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Contenu du stage"
$chart2.HasLegend = $true
$chart2.Legend.Position = -4152
$chart2.SeriesCollection(1).HasDataLabels = $true
$dLbls2 = $chart2.SeriesCollection(1).DataLabels()
$dLbls2.ShowValue = $true
$dLbls2.ShowPercentage = $true
$dLbls2.ShowCategoryName = $false
$dLbls2.ShowSeriesName = $false
$dLbls2.ShowLegendKey = $false

# --- Chart 3 : "Type du stage" (D25:E28, name in B25) ---------------------

$chartObj3 = $ws.ChartObjects().Add(205, 505, 300, 230)
$chart3 = $chartObj3.Chart
$chart3.ChartType = 5
$chart3.SetSourceData($ws.Range("D25:E28"))
$chart3.SeriesCollection(1).Name = "=Worksheet!`$B`$25"
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Type du stage"
$chart3.HasLegend = $true
$chart3.Legend.Position = -4152
$chart3.SeriesCollection(1).HasDataLabels = $true
$dLbls3 = $chart3.SeriesCollection(1).DataLabels()
$dLbls3.ShowValue = $true
$dLbls3.ShowPercentage = $true
$dLbls3.ShowCategoryName = $false
$dLbls3.ShowSeriesName = $false
$dLbls3.ShowLegendKey = $false

Write-Host "Charts created and labels repositioned."
